# LFM_Test_Plan.xlsx update:
# Mark a batch of Tier1-3 tests as "Passed" with a Last Run Date, and
# leave the Tiers1_3 sheet as the active/selected tab (with F22 selected).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Tiers1_3"
$ws2 = $wb.Worksheets.Item(2)   # "Tiers4_6"

# Rows on Tiers1_3 whose Status (col F) moves from "Planned" to "Passed"
# and which gain a Last Run Date (col G) of 10/23/2025 (serial 45953).
$rows = @(5, 11, 12, 13, 16, 17, 18, 19, 20, 21)

# Use an already-dated row as the formatting source so the new date cells
# pick up the existing short-date style instead of creating a new one.
$dateFormatSource = $ws1.Range("G2")
$dateFormatSource.Copy()

foreach ($r in $rows) {
    $ws1.Range("F$r").Value = "Passed"
    $gCell = $ws1.Range("G$r")
    $gCell.Value = 45953
    $gCell.PasteSpecial(-4122)  # xlPasteFormats
}

$excel.CutCopyMode = 0

# Make Tiers1_3 the active sheet/tab, with F22 selected (instead of Tiers4_6).
[void]$ws1.Activate()
[void]$ws1.Range("F22").Select()
